$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new entry -----------------------------------
# The list is kept in alphabetical order by name; "Respond to a lawsuit"
# belongs between row 14 ("Request time off work due to domestic
# violence") and what is currently row 15 ("Security deposit demand
# letter"). Insert a fresh row 15, which pushes the old rows 15-16 down
# to 16-17.
$ws.Rows.Item(15).Insert()

# --- 2. Fill in the new row --------------------------------------------
$ws.Range("A15").Value = "Respond to a lawsuit"
$ws.Range("B15").Value = "https://www.illinoislegalaid.org/legal-information/respond-lawsuit"

# --- 3. Rebuild the hyperlinks exactly in their final cell locations --
# Row insertion does not relocate the existing hyperlink anchors bound to
# the rows that moved (old B15/B16), so rebuild the hyperlink list from
# scratch in the correct, final row order.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"),  "https://www.illinoislegalaid.org/legal-information/appearance")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/fee-waiver")
$ws.Hyperlinks.Add($ws.Range("B5"),  "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter")
$ws.Hyperlinks.Add($ws.Range("B4"),  "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting")
$ws.Hyperlinks.Add($ws.Range("B9"),  "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter")
$ws.Hyperlinks.Add($ws.Range("B6"),  "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court")
$ws.Hyperlinks.Add($ws.Range("B7"),  "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court")
$ws.Hyperlinks.Add($ws.Range("B8"),  "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.illinoislegalaid.org/legal-information/respond-lawsuit")

# --- 4. Normalize formatting --------------------------------------------
# Hyperlinks.Add() stamps its own ad-hoc cell format; reassert the shared
# "Hyperlink" cell style on every hyperlinked cell in column B so they all
# reference the same style record (as they did before the edit).
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B15").Style = "Hyperlink"
$ws.Range("B16").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"
